# Edit script: applies the "3020 - IoT Vuln Research - STARTED" diff.
#
# Strategy:
#  - Small same-formatting run merges are done with a Find/Replace that
#    replaces text with itself (the engine naturally re-coalesces the
#    paragraph's runs when it rewrites the matched span).
#  - Paragraph-level "apply yellow highlight to this whole bullet" edits
#    (which need BOTH a run-level <w:rPr><w:highlight .../></w:rPr> AND a
#    paragraph-mark-level <w:pPr><w:rPr><w:highlight .../></w:rPr></w:pPr>)
#    are done with Range.InsertXML, which lets us write the exact target
#    OOXML for the paragraph (preserving its w14:paraId/rsid attributes).
#  - The new trailing space run after the hyperlink is added with
#    Range.InsertAfter at a point immediately before the paragraph mark.

$d = $word.ActiveDocument
$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParagraphByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) "... Must follow APA ..." bullet: merge the leading-space run and
#    the "Include sources..." run into a single run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(" Include sources and any personal experience.", $true, $false, $false, $false, $false, $true, 1, $false, " Include sources and any personal experience.", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) "Cover Page: ..." bullet -> whole-paragraph yellow highlight.
# ---------------------------------------------------------------------
$p = Get-ParagraphByText("Cover Page: Assignment name, course code, instructor name, your name, due date.")
$attrs = 'w14:paraId="06C70109" w14:textId="77777777" w:rsidR="00CC34ED" w:rsidRDefault="006752C1" w:rsidP="008723A1"'
$xml = '<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + $wordNs + ' ' + $attrs + '><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Cover Page: Assignment name, course code, instructor name, your name, due date.</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 3) "Table of Contents: Use Word's ..." bullet -> split the leading run
#    so "Table of Contents:" is its own (highlighted) run.
# ---------------------------------------------------------------------
$p = Get-ParagraphByText("Table of Contents: Use Word")
$attrs = 'w14:paraId="1B7F669A" w14:textId="669DDE48" w:rsidR="00CC34ED" w:rsidRDefault="006752C1" w:rsidP="008723A1"'
$xml = '<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + $wordNs + ' ' + $attrs + '><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>' +
       '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Table of Contents:</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> Use Word’s </w:t></w:r>' +
       '<w:r w:rsidR="00F358B1"><w:t xml:space="preserve">Table of Contents </w:t></w:r>' +
       '<w:r><w:t>feature.</w:t></w:r>' +
       '</w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 4) "For guidance on this tool ..." bullet -> append a trailing space
#    run right after the hyperlink (still inside the same paragraph).
# ---------------------------------------------------------------------
$p = Get-ParagraphByText("For guidance on this tool")
$r = $p.Range
$insertionPoint = $d.Range($r.End - 1, $r.End - 1)
$insertionPoint.InsertAfter(" ")

# ---------------------------------------------------------------------
# 5) "Headings: Apply Word styles ..." bullet -> whole-paragraph highlight.
# ---------------------------------------------------------------------
$p = Get-ParagraphByText("Headings: Apply Word styles")
$attrs = 'w14:paraId="018FDB7F" w14:textId="77777777" w:rsidR="00CC34ED" w:rsidRDefault="006752C1" w:rsidP="008723A1"'
$xml = '<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + $wordNs + ' ' + $attrs + '><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Headings: Apply Word styles (Heading 1, Heading 2).</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 6) "Header/Footer:" bullet -> whole-paragraph highlight.
# ---------------------------------------------------------------------
$p = Get-ParagraphByText("Header/Footer:")
$attrs = 'w14:paraId="06E4376F" w14:textId="77777777" w:rsidR="00CC34ED" w:rsidRDefault="006752C1" w:rsidP="008723A1"'
$xml = '<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + $wordNs + ' ' + $attrs + '><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Header/Footer:</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 7) "Header: Report title + your last name." bullet -> whole-paragraph
#    highlight.
# ---------------------------------------------------------------------
$p = Get-ParagraphByText("Header: Report title")
$attrs = 'w14:paraId="7C2F5BD0" w14:textId="5DC53AA6" w:rsidR="00CC34ED" w:rsidRDefault="006752C1" w:rsidP="008723A1"'
$xml = '<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + $wordNs + ' ' + $attrs + '><w:pPr><w:pStyle w:val="ListBullet"/><w:tabs><w:tab w:val="clear" w:pos="360"/><w:tab w:val="num" w:pos="720"/></w:tabs><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Header: Report title + your last name.</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 8) "Footer: Page numbers." bullet -> whole-paragraph highlight.
# ---------------------------------------------------------------------
$p = Get-ParagraphByText("Footer: Page numbers.")
$attrs = 'w14:paraId="3FE232A3" w14:textId="23F920B3" w:rsidR="00CC34ED" w:rsidRDefault="006752C1" w:rsidP="008723A1"'
$xml = '<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + $wordNs + ' ' + $attrs + '><w:pPr><w:pStyle w:val="ListBullet"/><w:tabs><w:tab w:val="num" w:pos="1897"/></w:tabs><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Footer: Page numbers.</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 9) "Font & Layout:" bullet -> whole-paragraph highlight.
# ---------------------------------------------------------------------
$p = Get-ParagraphByText("Font & Layout:")
$attrs = 'w14:paraId="63BD6104" w14:textId="77777777" w:rsidR="00CC34ED" w:rsidRDefault="006752C1" w:rsidP="008723A1"'
$xml = '<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + $wordNs + ' ' + $attrs + '><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Font &amp; Layout:</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 10) "Font: Arial, Calibri, or Aptos, size 12." bullet -> merge the
#     first two runs into one (highlighted), keep "." as its own
#     (highlighted) run, whole-paragraph highlight.
# ---------------------------------------------------------------------
$p = Get-ParagraphByText("Font: Arial, Calibri, or Aptos,")
$attrs = 'w14:paraId="733FD36E" w14:textId="524F7B77" w:rsidR="00CC34ED" w:rsidRDefault="006752C1" w:rsidP="008723A1"'
$xml = '<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + $wordNs + ' ' + $attrs + '><w:pPr><w:pStyle w:val="ListBullet"/><w:tabs><w:tab w:val="clear" w:pos="360"/><w:tab w:val="num" w:pos="720"/></w:tabs><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr>' +
       '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Font: Arial, Calibri, or Aptos, size 12</w:t></w:r>' +
       '<w:r w:rsidR="009C28E7"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>.</w:t></w:r>' +
       '</w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 11) "Margins: 1 inch; Line spacing: 1.5." bullet -> whole-paragraph
#     highlight.
# ---------------------------------------------------------------------
$p = Get-ParagraphByText("Margins: 1 inch; Line spacing: 1.5.")
$attrs = 'w14:paraId="729F2243" w14:textId="214A9E8E" w:rsidR="00CC34ED" w:rsidRDefault="006752C1" w:rsidP="008723A1"'
$xml = '<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + $wordNs + ' ' + $attrs + '><w:pPr><w:pStyle w:val="ListBullet"/><w:tabs><w:tab w:val="clear" w:pos="360"/><w:tab w:val="num" w:pos="720"/></w:tabs><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Margins: 1 inch; Line spacing: 1.5.</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

Write-Output "Edit complete."
